$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as plain text, preserving the workbook's existing (unstyled) look
function Set-TextValue($cellAddr, $value) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "27.964.51"
Set-TextValue "E2" "  +3.13%  "

# Row 3
Set-TextValue "D3" "1.571.25"
Set-TextValue "E3" "  +0.10%  "

# Row 4
Set-TextValue "D4" "0.979"
Set-TextValue "E4" "  -2.96%  "

# Row 5
Set-TextValue "D5" "209.68"
Set-TextValue "E5" "  -0.87%  "

# Row 6
Set-TextValue "D6" "0.494"
Set-TextValue "E6" "  +0.34%  "

# Row 7
Set-TextValue "D7" "0.981"
Set-TextValue "E7" "  -2.74%  "

# Row 8
Set-TextValue "D8" "23.41"
Set-TextValue "E8" "  +6.24%  "

# Row 9
Set-TextValue "D9" "0.249"
Set-TextValue "E9" "  +0.39%  "

# Row 10
Set-TextValue "D10" "0.0594"
Set-TextValue "E10" "  -0.61%  "

# Row 11
Set-TextValue "D11" "0.0868"
Set-TextValue "E11" "  +0.67%  "

# Row 12
Set-TextValue "D12" "1.802.23"
Set-TextValue "E12" "  +0.41%  "

# Row 13
Set-TextValue "D13" "1.581.13"
Set-TextValue "E13" "  +0.75%  "

# Row 14
Set-TextValue "D14" "3.73"
Set-TextValue "E14" "  -1.06%  "

# Row 15
Set-TextValue "D15" "0.520"
Set-TextValue "E15" "  -0.02%  "

# Row 16
Set-TextValue "D16" "28.166.07"
Set-TextValue "E16" "  +3.67%  "

# Row 17
Set-TextValue "D17" "62.84"
Set-TextValue "E17" "  +1.12%  "

# Row 18
Set-TextValue "D18" "234.32"
Set-TextValue "E18" "  +8.32%  "

# Row 19
Set-TextValue "D19" "7.50"
Set-TextValue "E19" "  +1.23%  "

# Row 20
Set-TextValue "D20" "0.0₃0702"
Set-TextValue "E20" "  -0.36%  "

# Row 21
Set-TextValue "D21" "0.971"
Set-TextValue "E21" "  -3.69%  "

# Row 22
Set-TextValue "D22" "4.10"
Set-TextValue "E22" "  -1.07%  "

# Row 23
Set-TextValue "D23" "9.36"
Set-TextValue "E23" "  +1.61%  "

# Row 24
Set-TextValue "D24" "1.93"
Set-TextValue "E24" "  -0.96%  "

# Row 25
Set-TextValue "D25" "148.59"
Set-TextValue "E25" "  -3.80%  "

# Row 26
Set-TextValue "D26" "15.31"
Set-TextValue "E26" "  +1.27%  "

# Row 27
Set-TextValue "D27" "0.107"
Set-TextValue "E27" "  +0.98%  "

# Row 28
Set-TextValue "D28" "6.55"
Set-TextValue "E28" "  -1.15%  "

# Row 29
Set-TextValue "D29" "0.977"
Set-TextValue "E29" "  -3.26%  "

# Row 30
Set-TextValue "D30" "1.13"
Set-TextValue "E30" "  -1.17%  "

# Row 31
Set-TextValue "D31" "0.0471"
Set-TextValue "E31" "  -0.55%  "

# Row 32
Set-TextValue "D32" "3.22"
Set-TextValue "E32" "  -0.59%  "

# Row 33
Set-TextValue "D33" "3.13"
Set-TextValue "E33" "  -2.13%  "

# Row 34
Set-TextValue "D34" "1.396.97"
Set-TextValue "E34" "  -2.68%  "

# Row 35
Set-TextValue "D35" "1.58"
Set-TextValue "E35" "  -1.71%  "

# Row 36
Set-TextValue "D36" "1.04"
Set-TextValue "E36" "  -6.11%  "

# Row 37
Set-TextValue "D37" "2.28"
Set-TextValue "E37" "  -2.93%  "

# Row 38
Set-TextValue "E38" "  +0.30%  "

# Row 39
Set-TextValue "D39" "0.546"
Set-TextValue "E39" "  +2.58%  "

# Row 40
Set-TextValue "D40" "2.43"

# Row 41
Set-TextValue "D41" "0.810"
Set-TextValue "E41" "  -0.07%  "

# Row 42
Set-TextValue "D42" "5.67"
Set-TextValue "E42" "  -3.59%  "

# Row 43
Set-TextValue "B43" "PaxDollar"
Set-TextValue "C43" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D43" "0.974"
Set-TextValue "E43" "  -3.48%  "

# Row 44
Set-TextValue "B44" "RenderToken"
Set-TextValue "C44" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D44" "1.85"
Set-TextValue "E44" "  +5.98%  "

# Row 45
Set-TextValue "E45" "  -3.50%  "

# Row 46
Set-TextValue "D46" "63.91"
Set-TextValue "E46" "  -1.05%  "

# Row 47
Set-TextValue "D47" "1.709.19"
Set-TextValue "E47" "  -0.11%  "

# Row 48
Set-TextValue "D48" "86.47"
Set-TextValue "E48" "  +0.77%  "

# Row 49
Set-TextValue "D49" "0.0523"
Set-TextValue "E49" "  +1.11%  "

# Row 50
Set-TextValue "D50" "40.79"
Set-TextValue "E50" "  +19.86%  "

# Row 51
Set-TextValue "D51" "0.0₇0983"
Set-TextValue "E51" "  -2.76%  "
